# Add a new "Greece" sheet (test data) right after "Croatia", based on the
# existing Croatia sheet layout, and make Greece the active/selected sheet.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")

# Duplicate the Croatia worksheet and place the copy immediately after it.
$croatia.Copy($null, $croatia)
$greece = $wb.Worksheets.Item("Croatia (2)")
$greece.Name = "Greece"

# Update the market name / ticket reference cells for the new country.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3164"

# Let Excel recompute the wrapped-text row heights for the rows whose
# content changed.
$greece.Rows.Item(3).EntireRow.AutoFit()
$greece.Rows.Item(4).EntireRow.AutoFit()

# Croatia is no longer the selected/active tab - reset its selection to a
# full-sheet selection like the other non-active sheets.
[void]$croatia.Range("A1:XFD1048576").Select()

# Make the new Greece sheet the active tab with its own selection.
[void]$greece.Activate()
[void]$greece.Range("D14").Select()
